$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (50P IOC Trip Pickup) for rows 2-20: 13 -> 2
for ($r = 2; $r -le 20; $r++) {
    $ws.Range("F$r").Value = 2
}

# Update column C (Max. Slip Voltage [%]) for the 480V relay rows: 5 -> 10
$rows480 = @(13, 14, 17, 18, 19, 20)
foreach ($r in $rows480) {
    $ws.Range("C$r").Value = 10
}

# Update the active selection on the sheet view
$ws.Range("F23").Select()
